$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell 'D2' '65.453.89'
Set-TextCell 'E2' '  +1.82%  '

# Row 3
Set-TextCell 'D3' '2.647.73'
Set-TextCell 'E3' '  +1.22%  '

# Row 4
Set-TextCell 'E4' '  +0.00%  '

# Row 5
Set-TextCell 'D5' '607.48'
Set-TextCell 'E5' '  +2.48%  '

# Row 6
Set-TextCell 'D6' '156.03'
Set-TextCell 'E6' '  +2.89%  '

# Row 7
Set-TextCell 'D7' '0.999'
Set-TextCell 'E7' '  -0.11%  '

# Row 8
Set-TextCell 'E8' '  -0.08%  '

# Row 9
Set-TextCell 'E9' '  +8.12%  '

# Row 10
Set-TextCell 'E10' '  +3.20%  '

# Row 11
Set-TextCell 'E11' '  +2.03%  '

# Row 12
Set-TextCell 'E12' '  +1.56%  '

# Row 13
Set-TextCell 'D13' '29.96'
Set-TextCell 'E13' '  +5.25%  '

# Row 14
Set-TextCell 'E14' '  +19.71%  '

# Row 15
Set-TextCell 'D15' '3.125.89'
Set-TextCell 'E15' '  +1.29%  '

# Row 16
Set-TextCell 'D16' '65.287.00'
Set-TextCell 'E16' '  +1.68%  '

# Row 17
Set-TextCell 'D17' '2.648.66'
Set-TextCell 'E17' '  -0.10%  '

# Row 18
Set-TextCell 'D18' '12.64'
Set-TextCell 'E18' '  +3.48%  '

# Row 19
Set-TextCell 'E19' '  +2.38%  '

# Row 20
Set-TextCell 'D20' '358.11'
Set-TextCell 'E20' '  +2.47%  '

# Row 21
Set-TextCell 'E21' '  +3.94%  '

# Row 22
Set-TextCell 'D22' '1.00'
Set-TextCell 'E22' '  -0.07%  '

# Row 23
Set-TextCell 'D23' '70.09'
Set-TextCell 'E23' '  +3.73%  '

# Row 24
Set-TextCell 'E24' '  +0.37%  '

# Row 25
Set-TextCell 'D25' '9.51'
Set-TextCell 'E25' '  +2.55%  '

# Row 26
Set-TextCell 'E26' '  +16.35%  '

# Row 27
Set-TextCell 'D27' '1.63'
Set-TextCell 'E27' '  -0.41%  '

# Row 28
Set-TextCell 'D28' '0.170'
Set-TextCell 'E28' '  +3.74%  '

# Row 29
Set-TextCell 'E29' '  -1.62%  '

# Row 30
Set-TextCell 'D30' '2.19'
Set-TextCell 'E30' '  +6.12%  '

# Row 31
Set-TextCell 'E31' '  +0.39%  '

# Row 32
Set-TextCell 'D32' '533.87'
Set-TextCell 'E32' '  -1.95%  '

# Row 33
Set-TextCell 'E33' '  -1.36%  '

# Row 34
Set-TextCell 'E34' '  -2.84%  '

# Row 35
Set-TextCell 'D35' '6.38'
Set-TextCell 'E35' '  +2.71%  '

# Row 36
Set-TextCell 'E36' '  +2.35%  '

# Row 37
Set-TextCell 'D37' '20.62'
Set-TextCell 'E37' '  +2.65%  '

# Row 38
Set-TextCell 'D38' '162.69'
Set-TextCell 'E38' '  -0.56%  '

# Row 39
Set-TextCell 'E39' '  +0.62%  '

# Row 40
Set-TextCell 'D40' '0.999'
Set-TextCell 'E40' '  -0.01%  '

# Row 41
Set-TextCell 'D41' '1.00'
Set-TextCell 'E41' '  +0.04%  '

# Row 42
Set-TextCell 'D42' '168.44'
Set-TextCell 'E42' '  +0.30%  '

# Row 43
Set-TextCell 'D43' '41.99'
Set-TextCell 'E43' '  +1.03%  '

# Row 44
Set-TextCell 'E44' '  +1.70%  '

# Row 45
Set-TextCell 'D45' '2.32'
Set-TextCell 'E45' '  +5.43%  '

# Row 46
Set-TextCell 'D46' '0.0610'
Set-TextCell 'E46' '  +2.21%  '

# Row 47
Set-TextCell 'E47' '  -0.72%  '

# Row 48
Set-TextCell 'B48' 'Mantle'
Set-TextCell 'C48' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 'D48' '0.654'
Set-TextCell 'E48' '  +2.49%  '

# Row 49
Set-TextCell 'B49' 'VeChain'
Set-TextCell 'C49' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D49' '0.0263'
Set-TextCell 'E49' '  +5.33%  '

# Row 50
Set-TextCell 'E50' '  -0.10%  '

# Row 51
Set-TextCell 'D51' '19.69'
Set-TextCell 'E51' '  +2.40%  '
